$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7544.077
$ws.Range("M69").Value = -10943.9999
$ws.Range("I69").Value = 3939.3333
$ws.Range("K69").Value = 11817.9999
$ws.Range("K72").Value = 35453.9997
$ws.Range("H72").Value = 7544.077
$ws.Range("M72").Value = -31085.9997
$ws.Range("I72").Value = 3939.3333
$ws.Range("H74").Value = 16834.334
$ws.Range("I74").Value = 16834.334
$ws.Range("M74").Value = -15898.334
$ws.Range("K74").Value = 16834.334
$ws.Range("N76").Value = -9686
$ws.Range("J76").Value = 9056
$ws.Range("K76").Value = 5972.8335
$ws.Range("H76").Value = 7000.5557
$ws.Range("I76").Value = 5972.8335
$ws.Range("M76").Value = -5657.8335
$ws.Range("L76").Value = 9056
$ws.Range("H77").Value = 16834.334
$ws.Range("K77").Value = 84171.67
$ws.Range("I77").Value = 16834.334
$ws.Range("M77").Value = -79491.67
$ws.Range("L79").Value = 9056
$ws.Range("I79").Value = 5972.8335
$ws.Range("H79").Value = 7000.5557
$ws.Range("N79").Value = -11240
$ws.Range("K79").Value = 5972.8335
$ws.Range("M79").Value = -4880.8335
$ws.Range("J79").Value = 9056
$ws.Range("N96").Value = -5748.4999
$ws.Range("L96").Value = 3002.4999
$ws.Range("M96").Value = -5433.25
$ws.Range("H96").Value = 1508
$ws.Range("I96").Value = 2268.75
$ws.Range("J96").Value = 1000.8333
$ws.Range("K96").Value = 6806.25
$ws.Range("H111").Value = 1365
$ws.Range("K111").Value = 4095
$ws.Range("M111").Value = -1028
$ws.Range("I111").Value = 1365
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K137").Value = 5812.7649
$ws.Range("H137").Value = 2687.3784
$ws.Range("M137").Value = -3262.7649
$ws.Range("I137").Value = 1937.5883

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -3715.2444
$ws.Range("I32").Value = 4002.2444
$ws.Range("K32").Value = 4002.2444
$ws.Range("H32").Value = 4895.7236
$ws.Range("I45").Value = 2509.111
$ws.Range("H45").Value = 3316.9375
$ws.Range("M45").Value = -2132.111
$ws.Range("K45").Value = 2509.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J110").Value = 35000
$ws.Range("H110").Value = 35000
$ws.Range("N110").Value = -43180
$ws.Range("L110").Value = 35000
$ws.Range("N122").Value = -80099.664
$ws.Range("H122").Value = 70299.664
$ws.Range("L122").Value = 70299.664
$ws.Range("J122").Value = 70299.664
$ws.Range("M134").Value = -3272.727000000001
$ws.Range("I134").Value = 1935.909
$ws.Range("H134").Value = 3044.85
$ws.Range("K134").Value = 5807.727000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J10").Value = 3335.6667
$ws.Range("M10").Value = -360.33334
$ws.Range("N10").Value = -3613.6667
$ws.Range("L10").Value = 3335.6667
$ws.Range("K10").Value = 499.33334
$ws.Range("H10").Value = 1917.5
$ws.Range("I10").Value = 499.33334
$ws.Range("L13").Value = 1092.3334
$ws.Range("J13").Value = 1092.3334
$ws.Range("N13").Value = -1370.3334
$ws.Range("H13").Value = 1092.3334
$ws.Range("N62").Value = -22250
$ws.Range("K62").Value = 3888.6667
$ws.Range("J62").Value = 21002
$ws.Range("M62").Value = -3264.6667
$ws.Range("L62").Value = 21002
$ws.Range("H62").Value = 9593.111000000001
$ws.Range("I62").Value = 3888.6667
$ws.Range("J65").Value = 21002
$ws.Range("L65").Value = 105010
$ws.Range("I65").Value = 3888.6667
$ws.Range("H65").Value = 9593.111000000001
$ws.Range("K65").Value = 19443.3335
$ws.Range("N65").Value = -111250
$ws.Range("M65").Value = -16323.3335
$ws.Range("N132").Value = -22702.334
$ws.Range("M132").Value = -4124.900000000001
$ws.Range("L132").Value = 17642.334
$ws.Range("J132").Value = 5880.778
$ws.Range("K132").Value = 6654.900000000001
$ws.Range("I132").Value = 2218.3
$ws.Range("H132").Value = 3063.487
$ws.Range("M134").Value = -1400.5386
$ws.Range("I134").Value = 1311.8462
$ws.Range("H134").Value = 2862.1052
$ws.Range("K134").Value = 3935.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M46").Value = -550.14284
$ws.Range("N46").Value = -9497
$ws.Range("L46").Value = 9315
$ws.Range("J46").Value = 3105
$ws.Range("H46").Value = 2141.238
$ws.Range("I46").Value = 213.71428
$ws.Range("K46").Value = 641.14284
$ws.Range("I63").Value = 0
$ws.Range("L63").Value = 39012.999
$ws.Range("J63").Value = 13004.333
$ws.Range("H63").Value = 13004.333
$ws.Range("N63").Value = -40510.999
$ws.Range("M63").ClearContents()
$ws.Range("K63").Value = 0
$ws.Range("H64").Value = 111119230
$ws.Range("J64").Value = 14760.5
$ws.Range("L64").Value = 44281.5
$ws.Range("N64").Value = -44821.5
$ws.Range("H66").Value = 13004.333
$ws.Range("M66").ClearContents()
$ws.Range("L66").Value = 117038.997
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("N66").Value = -124526.997
$ws.Range("J66").Value = 13004.333
$ws.Range("J67").Value = 14760.5
$ws.Range("H67").Value = 111119230
$ws.Range("L67").Value = 44281.5
$ws.Range("N67").Value = -46153.5
$ws.Range("K112").Value = 500021130
$ws.Range("H112").Value = 71433560
$ws.Range("J112").Value = 3436
$ws.Range("N112").Value = -12524
$ws.Range("M112").Value = -500020022
$ws.Range("L112").Value = 10308
$ws.Range("I112").Value = 166673710
$ws.Range("N122").Value = -35291.7148
$ws.Range("H122").Value = 2935.0588
$ws.Range("L122").Value = 30391.7148
$ws.Range("J122").Value = 3376.8572
$ws.Range("M126").Value = -3550
$ws.Range("H126").Value = 4088
$ws.Range("K126").Value = 8490
$ws.Range("I126").Value = 2830
$ws.Range("K137").Value = 2778
$ws.Range("N137").Value = -395775.75
$ws.Range("H137").Value = 79448.62
$ws.Range("J137").Value = 128525.25
$ws.Range("M137").Value = 2322
$ws.Range("L137").Value = 385575.75
$ws.Range("I137").Value = 926
$ws.Range("H140").Value = 2757.5386
$ws.Range("M140").Value = -39.66669999999976
$ws.Range("I140").Value = 1739.8889
$ws.Range("K140").Value = 5219.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L80").Value = 8001.5713
$ws.Range("N80").Value = -9997.5713
$ws.Range("M80").Value = -455780.62
$ws.Range("J80").Value = 8001.5713
$ws.Range("I80").Value = 456778.62
$ws.Range("H80").Value = 282254.22
$ws.Range("K80").Value = 456778.62
$ws.Range("J83").Value = 8001.5713
$ws.Range("N83").Value = -49991.85649999999
$ws.Range("I83").Value = 456778.62
$ws.Range("H83").Value = 282254.22
$ws.Range("L83").Value = 40007.85649999999
$ws.Range("K83").Value = 2283893.1
$ws.Range("M83").Value = -2278901.1
$ws.Range("N122").Value = -44001.499
$ws.Range("M122").Value = -77555.5
$ws.Range("I122").Value = 26668.5
$ws.Range("H122").Value = 19851.166
$ws.Range("K122").Value = 80005.5
$ws.Range("L122").Value = 39101.499
$ws.Range("J122").Value = 13033.833
$ws.Range("M126").Value = -3387.5
$ws.Range("H126").Value = 3290.1924
$ws.Range("L126").Value = 14552.5005
$ws.Range("J126").Value = 4850.8335
$ws.Range("N126").Value = -19492.5005
$ws.Range("K126").Value = 5857.5
$ws.Range("I126").Value = 1952.5
$ws.Range("M132").Value = -4861.1891
$ws.Range("K132").Value = 7391.1891
$ws.Range("I132").Value = 2463.7297
$ws.Range("H132").Value = 2829.2827
$ws.Range("L133").Value = 71995.91
$ws.Range("N133").Value = -82115.91
$ws.Range("J133").Value = 71995.91
$ws.Range("H133").Value = 71995.91
$ws.Range("H135").Value = 68355.69500000001
$ws.Range("L135").Value = 68355.69500000001
$ws.Range("J135").Value = 68355.69500000001
$ws.Range("N135").Value = -78495.69500000001
$ws.Range("N138").Value = -86512.664
$ws.Range("J138").Value = 76232.664
$ws.Range("H138").Value = 76232.664
$ws.Range("L138").Value = 76232.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 5465.1304
$ws.Range("M40").Value = -5329.1304
$ws.Range("K40").Value = 5465.1304
$ws.Range("H40").Value = 6596.2666
$ws.Range("N68").Value = -11042.818
$ws.Range("H68").Value = 6423.591
$ws.Range("L68").Value = 9544.817999999999
$ws.Range("J68").Value = 9544.817999999999
$ws.Range("N71").Value = -55212.09
$ws.Range("L71").Value = 47724.09
$ws.Range("J71").Value = 9544.817999999999
$ws.Range("H71").Value = 6423.591
$ws.Range("M132").Value = -7251.0386
$ws.Range("K132").Value = 9781.0386
$ws.Range("I132").Value = 3260.3462
$ws.Range("H132").Value = 3589.9556

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N46").Value = -63744.668
$ws.Range("L46").Value = 63282.668
$ws.Range("J46").Value = 63282.668
$ws.Range("H46").Value = 63282.668
$ws.Range("M122").Value = -4483.999899999999
$ws.Range("I122").Value = 2311.3333
$ws.Range("H122").Value = 3754.3845
$ws.Range("K122").Value = 6933.999899999999
$ws.Range("M132").Value = -3936.3329
$ws.Range("K132").Value = 6466.3329
$ws.Range("I132").Value = 2155.4443
$ws.Range("H132").Value = 2620.1765
$ws.Range("N134").Value = -194918.004
$ws.Range("L134").Value = 189848.004
$ws.Range("H134").Value = 63282.668
$ws.Range("J134").Value = 63282.668
